{"js": "// Full before->after text mapping for each of the 100 table cells, in\n// document (row-major) order. The source document has duplicate cell\n// texts (e.g. \"6+35=\" appears more than once), so replacements must be\n// applied cell-by-cell in order rather than as a single global\n// find/replace-all.\nconst pairs = [\n  [\"8+8=\", \"62-38=\"],\n  [\"93-24=\", \"29+35=\"],\n  [\"86+9=\", \"8+66=\"],\n  [\"98-9=\", \"61-12=\"],\n  [\"4+27=\", \"28+24=\"],\n  [\"65+28=\", \"67-39=\"],\n  [\"43-15=\", \"82-54=\"],\n  [\"47+44=\", \"56-8=\"],\n  [\"61-4=\", \"28+46=\"],\n  [\"19+4=\", \"68+9=\"],\n  [\"5+88=\", \"28+33=\"],\n  [\"6+35=\", \"23+9=\"],\n  [\"59+6=\", \"38+17=\"],\n  [\"75-58=\", \"17+68=\"],\n  [\"85-7=\", \"8+5=\"],\n  [\"40-23=\", \"4+69=\"],\n  [\"48+37=\", \"95-9=\"],\n  [\"91-56=\", \"67+28=\"],\n  [\"47+28=\", \"73-58=\"],\n  [\"49+36=\", \"29+9=\"],\n  [\"6+78=\", \"55-8=\"],\n  [\"62+9=\", \"28+5=\"],\n  [\"14+28=\", \"48+15=\"],\n  [\"70-23=\", \"36-18=\"],\n  [\"9+77=\", \"76-9=\"],\n  [\"28+4=\", \"74-7=\"],\n  [\"37+5=\", \"64-57=\"],\n  [\"90-15=\", \"91-2=\"],\n  [\"75-66=\", \"73-34=\"],\n  [\"39+24=\", \"72-43=\"],\n  [\"69+5=\", \"27+14=\"],\n  [\"31-8=\", \"84-56=\"],\n  [\"46-37=\", \"55+29=\"],\n  [\"34+58=\", \"12+69=\"],\n  [\"94-17=\", \"72-53=\"],\n  [\"29+49=\", \"28+38=\"],\n  [\"15+77=\", \"37+14=\"],\n  [\"17+35=\", \"80-48=\"],\n  [\"94-6=\", \"15+27=\"],\n  [\"10-2=\", \"13+38=\"],\n  [\"56-17=\", \"53+19=\"],\n  [\"41-14=\", \"91-23=\"],\n  [\"35+36=\", \"50-32=\"],\n  [\"36-29=\", \"35+58=\"],\n  [\"5+36=\", \"81-4=\"],\n  [\"53+18=\", \"69+6=\"],\n  [\"33+9=\", \"88+7=\"],\n  [\"46+37=\", \"90-29=\"],\n  [\"68-59=\", \"81-65=\"],\n  [\"9+74=\", \"9+29=\"],\n  [\"61-44=\", \"28+58=\"],\n  [\"7+54=\", \"26+36=\"],\n  [\"32+39=\", \"59+8=\"],\n  [\"26-18=\", \"13+8=\"],\n  [\"55+16=\", \"17+6=\"],\n  [\"8+23=\", \"44-18=\"],\n  [\"71-59=\", \"66-47=\"],\n  [\"57+26=\", \"35+18=\"],\n  [\"84-79=\", \"9+13=\"],\n  [\"86+7=\", \"49+16=\"],\n  [\"56+25=\", \"50-1=\"],\n  [\"29+33=\", \"26+67=\"],\n  [\"95-88=\", \"6+5=\"],\n  [\"18+38=\", \"21-18=\"],\n  [\"63-47=\", \"85-38=\"],\n  [\"81-34=\", \"14-9=\"],\n  [\"63-54=\", \"64+18=\"],\n  [\"84-79=\", \"94-85=\"],\n  [\"66+5=\", \"90-27=\"],\n  [\"7+74=\", \"42-15=\"],\n  [\"34-16=\", \"9+33=\"],\n  [\"58+19=\", \"31-6=\"],\n  [\"7+37=\", \"53-4=\"],\n  [\"66+6=\", \"18+39=\"],\n  [\"17+35=\", \"60-16=\"],\n  [\"25-17=\", \"38+3=\"],\n  [\"91-17=\", \"29+38=\"],\n  [\"47+39=\", \"55-48=\"],\n  [\"9+49=\", \"93-64=\"],\n  [\"13+48=\", \"66-48=\"],\n  [\"47+17=\", \"24-9=\"],\n  [\"91-7=\", \"18+74=\"],\n  [\"70-29=\", \"92-78=\"],\n  [\"66+27=\", \"25+38=\"],\n  [\"9+66=\", \"25+68=\"],\n  [\"34+28=\", \"94-78=\"],\n  [\"6+35=\", \"65+29=\"],\n  [\"84-5=\", \"80-46=\"],\n  [\"82-5=\", \"73-66=\"],\n  [\"58+24=\", \"91-24=\"],\n  [\"42-39=\", \"80-66=\"],\n  [\"53-9=\", \"32-15=\"],\n  [\"57+9=\", \"40-25=\"],\n  [\"42-34=\", \"49+49=\"],\n  [\"49+37=\", \"42-39=\"],\n  [\"27+37=\", \"29+7=\"],\n  [\"80-56=\", \"91-9=\"],\n  [\"95-7=\", \"33+49=\"],\n  [\"48+44=\", \"18+33=\"],\n  [\"49+2=\", \"19+32=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst cols = 5;\nfor (let i = 0; i < pairs.length; i++) {\n  const row = Math.floor(i / cols);\n  const col = i % cols;\n  const [oldText, newText] = pairs[i];\n  const cell = table.getCell(row, col);\n\n  // Search scoped to this single cell so duplicate values elsewhere in\n  // the table are not affected, and insertText(..., replace) on the\n  // found range swaps only the text while keeping the run's existing\n  // formatting (rFonts/sz) and the paragraph's properties intact.\n  const results = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`No match for \"${oldText}\" at row ${row}, col ${col}`);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Full before->after text mapping for each of the 100 table cells, in\n# document (row-major) order. The source document has duplicate cell\n# texts (e.g. \"6+35=\" appears more than once), so replacements must be\n# applied cell-by-cell (by position) in order rather than as a single\n# document-wide find/replace.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$cols = 5\n\n$pairs = @(\n  @(\"8+8=\", \"62-38=\"),\n  @(\"93-24=\", \"29+35=\"),\n  @(\"86+9=\", \"8+66=\"),\n  @(\"98-9=\", \"61-12=\"),\n  @(\"4+27=\", \"28+24=\"),\n  @(\"65+28=\", \"67-39=\"),\n  @(\"43-15=\", \"82-54=\"),\n  @(\"47+44=\", \"56-8=\"),\n  @(\"61-4=\", \"28+46=\"),\n  @(\"19+4=\", \"68+9=\"),\n  @(\"5+88=\", \"28+33=\"),\n  @(\"6+35=\", \"23+9=\"),\n  @(\"59+6=\", \"38+17=\"),\n  @(\"75-58=\", \"17+68=\"),\n  @(\"85-7=\", \"8+5=\"),\n  @(\"40-23=\", \"4+69=\"),\n  @(\"48+37=\", \"95-9=\"),\n  @(\"91-56=\", \"67+28=\"),\n  @(\"47+28=\", \"73-58=\"),\n  @(\"49+36=\", \"29+9=\"),\n  @(\"6+78=\", \"55-8=\"),\n  @(\"62+9=\", \"28+5=\"),\n  @(\"14+28=\", \"48+15=\"),\n  @(\"70-23=\", \"36-18=\"),\n  @(\"9+77=\", \"76-9=\"),\n  @(\"28+4=\", \"74-7=\"),\n  @(\"37+5=\", \"64-57=\"),\n  @(\"90-15=\", \"91-2=\"),\n  @(\"75-66=\", \"73-34=\"),\n  @(\"39+24=\", \"72-43=\"),\n  @(\"69+5=\", \"27+14=\"),\n  @(\"31-8=\", \"84-56=\"),\n  @(\"46-37=\", \"55+29=\"),\n  @(\"34+58=\", \"12+69=\"),\n  @(\"94-17=\", \"72-53=\"),\n  @(\"29+49=\", \"28+38=\"),\n  @(\"15+77=\", \"37+14=\"),\n  @(\"17+35=\", \"80-48=\"),\n  @(\"94-6=\", \"15+27=\"),\n  @(\"10-2=\", \"13+38=\"),\n  @(\"56-17=\", \"53+19=\"),\n  @(\"41-14=\", \"91-23=\"),\n  @(\"35+36=\", \"50-32=\"),\n  @(\"36-29=\", \"35+58=\"),\n  @(\"5+36=\", \"81-4=\"),\n  @(\"53+18=\", \"69+6=\"),\n  @(\"33+9=\", \"88+7=\"),\n  @(\"46+37=\", \"90-29=\"),\n  @(\"68-59=\", \"81-65=\"),\n  @(\"9+74=\", \"9+29=\"),\n  @(\"61-44=\", \"28+58=\"),\n  @(\"7+54=\", \"26+36=\"),\n  @(\"32+39=\", \"59+8=\"),\n  @(\"26-18=\", \"13+8=\"),\n  @(\"55+16=\", \"17+6=\"),\n  @(\"8+23=\", \"44-18=\"),\n  @(\"71-59=\", \"66-47=\"),\n  @(\"57+26=\", \"35+18=\"),\n  @(\"84-79=\", \"9+13=\"),\n  @(\"86+7=\", \"49+16=\"),\n  @(\"56+25=\", \"50-1=\"),\n  @(\"29+33=\", \"26+67=\"),\n  @(\"95-88=\", \"6+5=\"),\n  @(\"18+38=\", \"21-18=\"),\n  @(\"63-47=\", \"85-38=\"),\n  @(\"81-34=\", \"14-9=\"),\n  @(\"63-54=\", \"64+18=\"),\n  @(\"84-79=\", \"94-85=\"),\n  @(\"66+5=\", \"90-27=\"),\n  @(\"7+74=\", \"42-15=\"),\n  @(\"34-16=\", \"9+33=\"),\n  @(\"58+19=\", \"31-6=\"),\n  @(\"7+37=\", \"53-4=\"),\n  @(\"66+6=\", \"18+39=\"),\n  @(\"17+35=\", \"60-16=\"),\n  @(\"25-17=\", \"38+3=\"),\n  @(\"91-17=\", \"29+38=\"),\n  @(\"47+39=\", \"55-48=\"),\n  @(\"9+49=\", \"93-64=\"),\n  @(\"13+48=\", \"66-48=\"),\n  @(\"47+17=\", \"24-9=\"),\n  @(\"91-7=\", \"18+74=\"),\n  @(\"70-29=\", \"92-78=\"),\n  @(\"66+27=\", \"25+38=\"),\n  @(\"9+66=\", \"25+68=\"),\n  @(\"34+28=\", \"94-78=\"),\n  @(\"6+35=\", \"65+29=\"),\n  @(\"84-5=\", \"80-46=\"),\n  @(\"82-5=\", \"73-66=\"),\n  @(\"58+24=\", \"91-24=\"),\n  @(\"42-39=\", \"80-66=\"),\n  @(\"53-9=\", \"32-15=\"),\n  @(\"57+9=\", \"40-25=\"),\n  @(\"42-34=\", \"49+49=\"),\n  @(\"49+37=\", \"42-39=\"),\n  @(\"27+37=\", \"29+7=\"),\n  @(\"80-56=\", \"91-9=\"),\n  @(\"95-7=\", \"33+49=\"),\n  @(\"48+44=\", \"18+33=\"),\n  @(\"49+2=\", \"19+32=\")\n)\n\nfor ($i = 0; $i -lt $pairs.Length; $i++) {\n  $row = [int][Math]::Floor($i / $cols) + 1\n  $col = ($i % $cols) + 1\n  $oldText = $pairs[$i][0]\n  $newText = $pairs[$i][1]\n  $cell = $t.Cell($row, $col)\n  $r = $cell.Range\n\n  # Cell range text includes the trailing cell-mark (and a CR for multi-\n  # paragraph cells); strip those control characters before comparing so\n  # we only validate/replace the visible text, leaving the run formatting\n  # (rFonts/sz) and paragraph mark untouched.\n  $actual = $r.Text.TrimEnd([char]7).TrimEnd([char]13)\n  if ($actual -ne $oldText) {\n    throw \"Mismatch at row $row col ${col}: expected [$oldText] got [$actual]\"\n  }\n  $r.Text = $newText\n}\nWrite-Output \"Done: $($pairs.Length) replacements\"\n"}
